# Add 2022-Q4 data:
#  - Duplicate the "2022-Q3" sheet, place the duplicate right before it, rename
#    the duplicate to "2022-Q4" and refresh its figures with the new quarter's
#    numbers (the original "2022-Q3"/"2022-Q2" sheets stay untouched, just
#    shift right by one tab).
#  - Insert a matching summary row into "总计" for the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate "2022-Q3" -> new sheet placed before it -> rename to "2022-Q4"
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
# NOTE: after Copy(), the variable used to invoke it tracks the newly
# created (and now active) sheet rather than the original, so re-resolve
# the duplicate by its auto-generated name instead of trusting $q3/.Index.
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Refresh the new sheet's figures for 2022-Q4, keeping the same text-cell
# formatting (no numeric number-format baked in) as the source data.
$q4.Range("D2:G2").NumberFormat = "@"
$q4.Range("D2").Value = "4.76"
$q4.Range("E2").Value = "92.90"
$q4.Range("F2").Value = "9.22"
$q4.Range("G2").Value = "0.4389"
$q4.Range("D2:G2").Style = "Normal"

# ---------------------------------------------------------------------------
# 2. Insert the new quarter into the "总计" overview sheet
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Re-apply the data-row formatting (copied from the row right below, which
# holds the old data and already carries the correct style) to the new row.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$total.Range("B2:D2").Style = "Normal"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.44

# Renumber the sequential index column and refresh the values that shifted
# down so they read exactly like the new quarterly snapshot.
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.54

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.57

# Restore the original "last active tab" convention (the pre-edit workbook
# had its final sheet, 2022-Q2, as the selected tab).
$wb.Worksheets.Item("2022-Q2").Activate()
